$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3559
$ws.Range("I32").Value = 2999
$ws.Range("J32").Value = 3745.6667
$ws.Range("K32").Value = 2999
$ws.Range("L32").Value = 3745.6667
$ws.Range("M32").Value = -2673
$ws.Range("N32").Value = -4397.6667
# Row 58
$ws.Range("H58").Value = 1035.5714
$ws.Range("I58").Value = 1060.6666
$ws.Range("J58").Value = 1016.75
$ws.Range("K58").Value = 3181.9998
$ws.Range("L58").Value = 3050.25
$ws.Range("M58").Value = -3031.9998
$ws.Range("N58").Value = -3350.25
# Row 76
$ws.Range("H76").Value = 3875
$ws.Range("I76").Value = 3833.3333
$ws.Range("K76").Value = 3833.3333
$ws.Range("M76").Value = -3518.3333
# Row 79
$ws.Range("H79").Value = 3875
$ws.Range("I79").Value = 3833.3333
$ws.Range("K79").Value = 3833.3333
$ws.Range("M79").Value = -2741.3333
# Row 92
$ws.Range("H92").Value = 1423.8
$ws.Range("I92").Value = 257.84616
$ws.Range("J92").Value = 9002.5
$ws.Range("K92").Value = 257.84616
$ws.Range("L92").Value = 9002.5
$ws.Range("M92").Value = 990.1538399999999
$ws.Range("N92").Value = -11498.5
# Row 134
$ws.Range("H134").Value = 100000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 100000
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -110140
# Row 135
$ws.Range("H135").Value = 2152.1177
$ws.Range("I135").Value = 772.4
$ws.Range("K135").Value = 6951.599999999999
$ws.Range("M135").Value = -4416.599999999999
# Row 136
$ws.Range("H136").Value = 125000
$ws.Range("J136").Value = 125000
$ws.Range("L136").Value = 125000
$ws.Range("N136").Value = -135200
# Row 137
$ws.Range("H137").Value = 79421380
$ws.Range("J137").Value = 2952357.2
$ws.Range("L137").Value = 8857071.600000001
$ws.Range("N137").Value = -8862171.600000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2817.7144
$ws.Range("I2").Value = 1623.9445
$ws.Range("J2").Value = 9980.333000000001
$ws.Range("K2").Value = 1623.9445
$ws.Range("L2").Value = 9980.333000000001
$ws.Range("M2").Value = -1510.9445
$ws.Range("N2").Value = -10206.333
# Row 74
$ws.Range("H74").Value = 2909362.8
$ws.Range("J74").Value = 5991.1816
$ws.Range("L74").Value = 5991.1816
$ws.Range("N74").Value = -7739.1816
# Row 77
$ws.Range("H77").Value = 2909362.8
$ws.Range("J77").Value = 5991.1816
$ws.Range("L77").Value = 29955.908
$ws.Range("N77").Value = -38691.908
# Row 116
$ws.Range("H116").Value = 2817.7144
$ws.Range("I116").Value = 1623.9445
$ws.Range("J116").Value = 9980.333000000001
$ws.Range("K116").Value = 1623.9445
$ws.Range("L116").Value = 9980.333000000001
$ws.Range("M116").Value = 670.0554999999999
$ws.Range("N116").Value = -14568.333
# Row 122
$ws.Range("H122").Value = 3394.5
$ws.Range("I122").Value = 3297.5
$ws.Range("K122").Value = 9892.5
$ws.Range("M122").Value = -7442.5
# Row 132
$ws.Range("H132").Value = 743513
$ws.Range("I132").Value = 880078.6
$ws.Range("J132").Value = 6058.8
$ws.Range("K132").Value = 2640235.8
$ws.Range("L132").Value = 18176.4
$ws.Range("M132").Value = -2637705.8
$ws.Range("N132").Value = -23236.4

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2817.7144
$ws.Range("I3").Value = 1623.9445
$ws.Range("J3").Value = 9980.333000000001
$ws.Range("K3").Value = 1623.9445
$ws.Range("L3").Value = 9980.333000000001
$ws.Range("M3").Value = -1509.9445
$ws.Range("N3").Value = -10208.333
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 86
$ws.Range("H86").Value = 1800.826
$ws.Range("I86").Value = 1724.6428
$ws.Range("J86").Value = 1919.3334
$ws.Range("K86").Value = 1724.6428
$ws.Range("L86").Value = 1919.3334
$ws.Range("M86").Value = -601.6428000000001
$ws.Range("N86").Value = -4165.3334
# Row 89
$ws.Range("H89").Value = 1800.826
$ws.Range("I89").Value = 1724.6428
$ws.Range("J89").Value = 1919.3334
$ws.Range("K89").Value = 8623.214
$ws.Range("L89").Value = 9596.666999999999
$ws.Range("M89").Value = -3007.214
$ws.Range("N89").Value = -20828.667
# Row 99
$ws.Range("H99").Value = 15659.6
$ws.Range("I99").Value = 6667.7144
$ws.Range("J99").Value = 36640.668
$ws.Range("K99").Value = 6667.7144
$ws.Range("L99").Value = 36640.668
$ws.Range("M99").Value = -5169.7144
$ws.Range("N99").Value = -39636.668
# Row 134
$ws.Range("H134").Value = 849757.6
$ws.Range("I134").Value = 1273675
$ws.Range("J134").Value = 312795.66
$ws.Range("K134").Value = 3821025
$ws.Range("L134").Value = 938386.98
$ws.Range("M134").Value = -3818490
$ws.Range("N134").Value = -943456.98
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 76
$ws.Range("H76").Value = 5932.3335
$ws.Range("I76").Value = 5932.3335
$ws.Range("K76").Value = 5932.3335
$ws.Range("M76").Value = -5617.3335
# Row 79
$ws.Range("H79").Value = 5932.3335
$ws.Range("I79").Value = 5932.3335
$ws.Range("K79").Value = 5932.3335
$ws.Range("M79").Value = -4840.3335

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 142.3
$ws.Range("I40").Value = 67.5
$ws.Range("J40").Value = 316.83334
$ws.Range("K40").Value = 270
$ws.Range("L40").Value = 1267.33336
$ws.Range("M40").Value = -201
$ws.Range("N40").Value = -1405.33336
# Row 75
$ws.Range("H75").Value = 4947.7334
$ws.Range("J75").Value = 7250.3335
$ws.Range("L75").Value = 21751.0005
$ws.Range("N75").Value = -23747.0005
# Row 78
$ws.Range("H78").Value = 4947.7334
$ws.Range("J78").Value = 7250.3335
$ws.Range("L78").Value = 65253.0015
$ws.Range("N78").Value = -75237.0015

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 7376.7
$ws.Range("I9").Value = 1467.25
$ws.Range("K9").Value = 1467.25
$ws.Range("M9").Value = -1297.25
# Row 122
$ws.Range("H122").Value = 50268.09
$ws.Range("I122").Value = 86853.086
$ws.Range("J122").Value = 6366.1
$ws.Range("K122").Value = 260559.258
$ws.Range("L122").Value = 19098.3
$ws.Range("M122").Value = -258109.258
$ws.Range("N122").Value = -23998.3
# Row 132
$ws.Range("H132").Value = 248360.1
$ws.Range("I132").Value = 315867.78
$ws.Range("J132").Value = 9014.637000000001
$ws.Range("K132").Value = 947603.3400000001
$ws.Range("L132").Value = 27043.911
$ws.Range("M132").Value = -945073.3400000001
$ws.Range("N132").Value = -32103.911

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 75
$ws.Range("H75").Value = 48500
# Row 78
$ws.Range("H78").Value = 48500
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
# Row 120
$ws.Range("H120").Value = 98558.39999999999
$ws.Range("J120").Value = 98558.39999999999
$ws.Range("L120").Value = 98558.39999999999
$ws.Range("N120").Value = -108234.4
# Row 134
$ws.Range("H134").Value = 45866.668
$ws.Range("J134").Value = 45866.668
$ws.Range("L134").Value = 45866.668
$ws.Range("N134").Value = -56006.668
# Row 135
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1805.4
$ws.Range("I100").Value = 1548.1875
$ws.Range("K100").Value = 3096.375
$ws.Range("M100").Value = -2555.375
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
# Row 136
$ws.Range("H136").Value = 33984916
$ws.Range("I136").Value = 35666896
$ws.Range("K136").Value = 107000688
$ws.Range("M136").Value = -106998138
